# Auto commit at 2025-10-23  7:22:58.92
# Refresh the monthly "Metrics" figures (B2:B13) with the latest totals.
# The "today" sheet pulls these via formulas (=Metrics!Bn), so its B/E/F
# columns recompute automatically; TODAY()-1 in today!A1 also re-evaluates.

$wb = $excel.ActiveWorkbook

$metrics = $wb.Worksheets.Item("Metrics")

$metrics.Range("B2").Value2  = 304944.02
$metrics.Range("B3").Value2  = 249913.44
$metrics.Range("B4").Value2  = 97161.819999999978
$metrics.Range("B5").Value2  = 12153
$metrics.Range("B6").Value2  = 4672075.49
$metrics.Range("B7").Value2  = 3939732.1099999994
$metrics.Range("B8").Value2  = 1367763.96
$metrics.Range("B9").Value2  = 181154
$metrics.Range("B10").Value2 = 33137399.290999822
$metrics.Range("B11").Value2 = 31214953.629999999
$metrics.Range("B12").Value2 = 11649472.850000001
$metrics.Range("B13").Value2 = 1278781

# Match the recorded selection on the Metrics sheet.
$metrics.Activate()
$metrics.Range("D15").Select()

# Recalculate so dependent sheets (today!) pick up the new Metrics values.
$excel.Calculate()

# Restore focus/selection to the "today" sheet (the tab active in the file).
$today = $wb.Worksheets.Item("today")
$today.Activate()
$today.Range("E7").Select()
